$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Select column D (the "Credits Enrolled" column) and delete it entirely,
# shifting Status (old E) and the following columns left.
$ws.Columns.Item(4).Select() | Out-Null
$ws.Columns.Item(4).Delete() | Out-Null
